$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.6363636363636364
$ws.Range("C2").Value = 0.7368421052631579
$ws.Range("D2").Value = 0.6829268292682926
$ws.Range("B3").Value = 0.7959183673469388
$ws.Range("C3").Value = 0.7090909090909091
$ws.Range("D3").Value = 0.75
$ws.Range("B4").Value = 0.7204301075268817
$ws.Range("C4").Value = 0.7204301075268817
$ws.Range("D4").Value = 0.7204301075268817
$ws.Range("E4").Value = 0.7204301075268817
$ws.Range("B5").Value = 0.7161410018552876
$ws.Range("C5").Value = 0.7229665071770335
$ws.Range("D5").Value = 0.7164634146341463
$ws.Range("B6").Value = 0.7307239611387076
$ws.Range("C6").Value = 0.7204301075268817
$ws.Range("D6").Value = 0.7225937581956465
$ws.Range("B7").Value = 0.5909090909090909
$ws.Range("C7").Value = 0.6842105263157895
$ws.Range("D7").Value = 0.6341463414634148
$ws.Range("B8").Value = 0.7551020408163265
$ws.Range("D8").Value = 0.7115384615384616
$ws.Range("B9").Value = 0.6774193548387096
$ws.Range("C9").Value = 0.6774193548387096
$ws.Range("D9").Value = 0.6774193548387096
$ws.Range("E9").Value = 0.6774193548387096
$ws.Range("B10").Value = 0.6730055658627088
$ws.Range("C10").Value = 0.6784688995215311
$ws.Range("D10").Value = 0.6728424015009382
$ws.Range("B11").Value = 0.6880124483811121
$ws.Range("C11").Value = 0.6774193548387096
$ws.Range("D11").Value = 0.6799158748411306
$ws.Range("B22").Value = 0.5961538461538461
$ws.Range("C22").Value = 0.8157894736842105
$ws.Range("D22").Value = 0.6888888888888889
$ws.Range("B23").Value = 0.8292682926829268
$ws.Range("C23").Value = 0.6181818181818182
$ws.Range("D23").Value = 0.7083333333333334
$ws.Range("B24").Value = 0.6989247311827957
$ws.Range("C24").Value = 0.6989247311827957
$ws.Range("D24").Value = 0.6989247311827957
$ws.Range("E24").Value = 0.6989247311827957
$ws.Range("B25").Value = 0.7127110694183865
$ws.Range("C25").Value = 0.7169856459330144
$ws.Range("D25").Value = 0.6986111111111111
$ws.Range("B26").Value = 0.734017228509754
$ws.Range("C26").Value = 0.6989247311827957
$ws.Range("D26").Value = 0.7003882915173238
